$d = $word.ActiveDocument

# The document's single section has a "first page" header/footer pair plus
# a "default" (other pages) header/footer pair:
#   Headers.Item(2) / Footers.Item(2) -> wdHeaderFooterFirstPage
#   Footers.Item(1)                   -> wdHeaderFooterPrimary
# Three pictures inside those headers/footers need their display/default
# names swapped, matching how Word's picture-name counters were
# renumbered on the source machine:
#   - BTEC logo (first-page header)      image1.jpg -> image2.jpg
#   - Pearson logo (first-page footer)   image2.png -> image1.png
#   - Pearson logo (default/other footer) image2.png -> image1.png

$sec = $d.Sections.Item(1)

$headerFirst = $sec.Headers.Item(2)
$btecLogo = $headerFirst.Range.InlineShapes.Item(1)
$btecLogo.Name = "image2.jpg"

$footerDefault = $sec.Footers.Item(1)
$pearsonLogoDefault = $footerDefault.Range.InlineShapes.Item(1)
$pearsonLogoDefault.Name = "image1.png"

$footerFirst = $sec.Footers.Item(2)
$pearsonLogoFirst = $footerFirst.Range.InlineShapes.Item(1)
$pearsonLogoFirst.Name = "image1.png"
